# Natmi following Dr Hou advice
# Update LR-pair values for rows 2-10 and append new rows 11-16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Areg"
$ws.Range("C2").Value = "Egfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.484834666666667
$ws.Range("H2").Value = 4.454504
$ws.Range("I2").Value = 0.4307162850350085
$ws.Range("J2").Value = 0.4307162850350084
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.57413
$ws.Range("N2").Value = 4.72239
$ws.Range("O2").Value = 0.02024862668342525
$ws.Range("P2").Value = 0.02024862668342525
$ws.Range("Q2").Value = 2.33732279384
$ws.Range("R2").Value = 21.03590514456
$ws.Range("S2").Value = 0.008721413262145667
$ws.Range("T2").Value = 0.008721413262145665

# Row 3
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Areg"
$ws.Range("C3").Value = "Egfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.484834666666667
$ws.Range("H3").Value = 4.454504
$ws.Range("I3").Value = 0.4307162850350085
$ws.Range("J3").Value = 0.4307162850350084
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 62.503947
$ws.Range("N3").Value = 187.511841
$ws.Range("O3").Value = 0.804011796385049
$ws.Range("P3").Value = 0.8040117963850492
$ws.Range("Q3").Value = 92.808027309096
$ws.Range("R3").Value = 835.2722457818641
$ws.Range("S3").Value = 0.346300974063292
$ws.Range("T3").Value = 0.346300974063292

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Areg"
$ws.Range("C4").Value = "Egfr"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.484834666666667
$ws.Range("H4").Value = 4.454504
$ws.Range("I4").Value = 0.4307162850350085
$ws.Range("J4").Value = 0.4307162850350084
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.182903
$ws.Range("N4").Value = 0.548709
$ws.Range("O4").Value = 0.002352750132631058
$ws.Range("P4").Value = 0.002352750132631058
$ws.Range("Q4").Value = 0.2715807150373334
$ws.Range("R4").Value = 2.444226435336
$ws.Range("S4").Value = 0.001013367796742473
$ws.Range("T4").Value = 0.001013367796742473

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Areg"
$ws.Range("C5").Value = "Egfr"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.484834666666667
$ws.Range("H5").Value = 4.454504
$ws.Range("I5").Value = 0.4307162850350085
$ws.Range("J5").Value = 0.4307162850350084
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.09380766666666666
$ws.Range("N5").Value = 0.281423
$ws.Range("O5").Value = 0.001206683324996365
$ws.Range("P5").Value = 0.001206683324996365
$ws.Range("Q5").Value = 0.1392888754657778
$ws.Range("R5").Value = 1.253599879192
$ws.Range("S5").Value = 0.000519738158956126
$ws.Range("T5").Value = 0.0005197381589561259

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Areg"
$ws.Range("C6").Value = "Egfr"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.484834666666667
$ws.Range("H6").Value = 4.454504
$ws.Range("I6").Value = 0.4307162850350085
$ws.Range("J6").Value = 0.4307162850350084
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.38529933333333
$ws.Range("N6").Value = 40.155898
$ws.Range("O6").Value = 0.1721801434738983
$ws.Range("P6").Value = 0.1721801434738983
$ws.Range("Q6").Value = 19.87495647384356
$ws.Range("R6").Value = 178.874608264592
$ws.Range("S6").Value = 0.07416079175387222
$ws.Range("T6").Value = 0.0741607917538722

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Areg"
$ws.Range("C7").Value = "Egfr"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.146843333333333
$ws.Range("H7").Value = 3.44053
$ws.Range("I7").Value = 0.332672795927784
$ws.Range("J7").Value = 0.332672795927784
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.57413
$ws.Range("N7").Value = 4.72239
$ws.Range("O7").Value = 0.02024862668342525
$ws.Range("P7").Value = 0.02024862668342525
$ws.Range("Q7").Value = 1.8052804963
$ws.Range("R7").Value = 16.2475244667
$ws.Range("S7").Value = 0.006736167252473008
$ws.Range("T7").Value = 0.006736167252473008

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Areg"
$ws.Range("C8").Value = "Egfr"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.146843333333333
$ws.Range("H8").Value = 3.44053
$ws.Range("I8").Value = 0.332672795927784
$ws.Range("J8").Value = 0.332672795927784
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 62.503947
$ws.Range("N8").Value = 187.511841
$ws.Range("O8").Value = 0.804011796385049
$ws.Range("P8").Value = 0.8040117963850492
$ws.Range("Q8").Value = 71.68223492397
$ws.Range("R8").Value = 645.14011431573
$ws.Range("S8").Value = 0.2674728522623344
$ws.Range("T8").Value = 0.2674728522623345

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Areg"
$ws.Range("C9").Value = "Egfr"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.146843333333333
$ws.Range("H9").Value = 3.44053
$ws.Range("I9").Value = 0.332672795927784
$ws.Range("J9").Value = 0.332672795927784
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.182903
$ws.Range("N9").Value = 0.548709
$ws.Range("O9").Value = 0.002352750132631058
$ws.Range("P9").Value = 0.002352750132631058
$ws.Range("Q9").Value = 0.2097610861966667
$ws.Range("R9").Value = 1.88784977577
$ws.Range("S9").Value = 0.0007826959647418387
$ws.Range("T9").Value = 0.0007826959647418387

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Areg"
$ws.Range("C10").Value = "Egfr"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.146843333333333
$ws.Range("H10").Value = 3.44053
$ws.Range("I10").Value = 0.332672795927784
$ws.Range("J10").Value = 0.332672795927784
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.09380766666666666
$ws.Range("N10").Value = 0.281423
$ws.Range("O10").Value = 0.001206683324996365
$ws.Range("P10").Value = 0.001206683324996365
$ws.Range("Q10").Value = 0.1075826971322222
$ws.Range("R10").Value = 0.9682442741899999
$ws.Range("S10").Value = 0.0004014307155259755
$ws.Range("T10").Value = 0.0004014307155259755

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Areg"
$ws.Range("C11").Value = "Egfr"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.146843333333333
$ws.Range("H11").Value = 3.44053
$ws.Range("I11").Value = 0.332672795927784
$ws.Range("J11").Value = 0.332672795927784
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 13.38529933333333
$ws.Range("N11").Value = 40.155898
$ws.Range("O11").Value = 0.1721801434738983
$ws.Range("P11").Value = 0.1721801434738983
$ws.Range("Q11").Value = 15.35084130510444
$ws.Range("R11").Value = 138.15757174594
$ws.Range("S11").Value = 0.05727964973270872
$ws.Range("T11").Value = 0.05727964973270872

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Areg"
$ws.Range("C12").Value = "Egfr"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.8156833333333333
$ws.Range("H12").Value = 2.44705
$ws.Range("I12").Value = 0.2366109190372076
$ws.Range("J12").Value = 0.2366109190372076
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.57413
$ws.Range("N12").Value = 4.72239
$ws.Range("O12").Value = 0.02024862668342525
$ws.Range("P12").Value = 0.02024862668342525
$ws.Range("Q12").Value = 1.2839916055
$ws.Range("R12").Value = 11.5559244495
$ws.Range("S12").Value = 0.004791046168806571
$ws.Range("T12").Value = 0.004791046168806571

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Areg"
$ws.Range("C13").Value = "Egfr"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.8156833333333333
$ws.Range("H13").Value = 2.44705
$ws.Range("I13").Value = 0.2366109190372076
$ws.Range("J13").Value = 0.2366109190372076
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 62.503947
$ws.Range("N13").Value = 187.511841
$ws.Range("O13").Value = 0.804011796385049
$ws.Range("P13").Value = 0.8040117963850492
$ws.Range("Q13").Value = 50.98342783545
$ws.Range("R13").Value = 458.85085051905
$ws.Range("S13").Value = 0.1902379700594227
$ws.Range("T13").Value = 0.1902379700594227

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Areg"
$ws.Range("C14").Value = "Egfr"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.8156833333333333
$ws.Range("H14").Value = 2.44705
$ws.Range("I14").Value = 0.2366109190372076
$ws.Range("J14").Value = 0.2366109190372076
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.182903
$ws.Range("N14").Value = 0.548709
$ws.Range("O14").Value = 0.002352750132631058
$ws.Range("P14").Value = 0.002352750132631058
$ws.Range("Q14").Value = 0.1491909287166667
$ws.Range("R14").Value = 1.34271835845
$ws.Range("S14").Value = 0.0005566863711467467
$ws.Range("T14").Value = 0.0005566863711467467

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Areg"
$ws.Range("C15").Value = "Egfr"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8156833333333333
$ws.Range("H15").Value = 2.44705
$ws.Range("I15").Value = 0.2366109190372076
$ws.Range("J15").Value = 0.2366109190372076
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.09380766666666666
$ws.Range("N15").Value = 0.281423
$ws.Range("O15").Value = 0.001206683324996365
$ws.Range("P15").Value = 0.001206683324996365
$ws.Range("Q15").Value = 0.07651735023888888
$ws.Range("R15").Value = 0.68865615215
$ws.Range("S15").Value = 0.0002855144505142633
$ws.Range("T15").Value = 0.0002855144505142633

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Areg"
$ws.Range("C16").Value = "Egfr"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.8156833333333333
$ws.Range("H16").Value = 2.44705
$ws.Range("I16").Value = 0.2366109190372076
$ws.Range("J16").Value = 0.2366109190372076
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 13.38529933333333
$ws.Range("N16").Value = 40.155898
$ws.Range("O16").Value = 0.1721801434738983
$ws.Range("P16").Value = 0.1721801434738983
$ws.Range("Q16").Value = 10.91816557787778
$ws.Range("R16").Value = 98.2634902009
$ws.Range("S16").Value = 0.04073970198731733
$ws.Range("T16").Value = 0.04073970198731733
